# Students_Details.xlsx edit:
# Insert a new "Current Address" column before the existing "State" column
# (old column K), pushing State/City/Age/Salary/Department one column to
# the right, and fill in the four students' current-address values.
# Also nudge the selection / column widths to match the saved state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at K (old K = "State" shifts to L, etc.)
$ws.Columns.Item(11).Insert()

# New column header + data (literal backslash-n sequences, not real newlines)
$ws.Range("K1").Value = "Current Address"
$ws.Range("K2").Value = '301 Town 481, Brookville\n kansas 42129\n United States'
$ws.Range("K3").Value = '302 Town 481, Brookville\n kansas 42129\n United States'
$ws.Range("K4").Value = '303 Town 481, Brookville\n kansas 42129\n United States'
$ws.Range("K5").Value = '304 Town 481, Brookville\n kansas 42129\n United States'

# Column width tweaks (Mobile column narrows, new Current Address column
# gets the same wide "address" style width as the existing Address column)
$ws.Columns.Item(5).ColumnWidth = 10.166666666666666
$ws.Columns.Item(11).ColumnWidth = 49.666666666666664

# Update the remembered selection
$ws.Range("K11").Select() | Out-Null
